$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.847.80"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.438.30"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'560.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'162.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "'0.168"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.98%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  -5.66%  "
$ws.Range("D13").Value = "'0.0000176"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("D14").Value = "68.713.57"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "2.886.38"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "'23.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "2.440.35"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'339.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'3.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.71%  "
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'67.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").Value = "2.567.52"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "'7.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'428.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'158.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'18.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +3.97%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'1.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").Value = "'130.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'0.0924"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +1.14%  "
